$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append below the existing data (rows 321-328),
# corresponding to dates 2021-02-21 .. 2021-02-28.
$newRows = @(
    @{ Row = 321; Date = "2021-02-21"; Scheduled = 50; Tracked = 48 },
    @{ Row = 322; Date = "2021-02-22"; Scheduled = 61; Tracked = 60 },
    @{ Row = 323; Date = "2021-02-23"; Scheduled = 52; Tracked = 49 },
    @{ Row = 324; Date = "2021-02-24"; Scheduled = 52; Tracked = 47 },
    @{ Row = 325; Date = "2021-02-25"; Scheduled = 53; Tracked = 52 },
    @{ Row = 326; Date = "2021-02-26"; Scheduled = 49; Tracked = 49 },
    @{ Row = 327; Date = "2021-02-27"; Scheduled = 45; Tracked = 45 },
    @{ Row = 328; Date = "2021-02-28"; Scheduled = 44; Tracked = 44 }
)

foreach ($item in $newRows) {
    $row = $item.Row
    $prevRow = $row - 1

    # Copy the row above first so the new row inherits the same cell
    # styles/number formats (text style for col A, integer style for
    # B/C, percentage style for col D) as the rest of the table.
    $ws.Range("A" + $prevRow + ":D" + $prevRow).Copy($ws.Range("A" + $row + ":D" + $row))

    # Use Formula (not Value) for the date text so it stays a plain text
    # string like the rest of column A, rather than being auto-converted
    # to a date serial number.
    $ws.Range("A" + $row).Formula = $item.Date
    $ws.Range("B" + $row).Formula = $item.Scheduled
    $ws.Range("C" + $row).Formula = $item.Tracked
    $ws.Range("D" + $row).Formula = "=C" + $row + "/B" + $row
}

# Reflect the "select all" / scrolled view captured in the saved sheet view.
$ws.Cells.Select()
